# Timetracking workbook: add new entries "in advance" for the final days
# of the project (27.01.2025 - 31.01.2025), as described in the commit
# message "Added Timetracking in advance in case I break everything".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zeitaufzeichnung")

# Row 78: 27.01.2025 - Anpassen der Dokumentation (6h)
$ws.Cells.Item(78, 1).Value = "27.01.2025"
$ws.Cells.Item(78, 2).Value = 6
$ws.Cells.Item(78, 3).Value = "Anpassen der Dokumentation"

# Row 79: 28.01.2025 - Code optimierung (6h)
$ws.Cells.Item(79, 1).Value = "28.01.2025"
$ws.Cells.Item(79, 2).Value = 6
$ws.Cells.Item(79, 3).Value = "Code optimierung"

# Row 80: 29.01.2025 - Code optimierung (5h)
$ws.Cells.Item(80, 1).Value = "29.01.2025"
$ws.Cells.Item(80, 2).Value = 5
$ws.Cells.Item(80, 3).Value = "Code optimierung"

# Row 81: 30.01.2025 - Cleanup des git repos und Tests (5h)
$ws.Cells.Item(81, 1).Value = "30.01.2025"
$ws.Cells.Item(81, 2).Value = 5
$ws.Cells.Item(81, 3).Value = "Cleanup des git repos und Tests"

# Row 82: 31.01.2025 - Tests und Abgabe des Projekts (1h)
$ws.Cells.Item(82, 1).Value = "31.01.2025"
$ws.Cells.Item(82, 2).Value = 1
$ws.Cells.Item(82, 3).Value = "Tests und Abgabe des Projekts"

# The "total" row (96) already sums B2:B95 via formula, so it will
# automatically pick up the newly added hours once recalculated.

# Scroll the sheet view up by one row to match the author's saved view
$excel.Goto($ws.Range("A72"), $true)
